# Rename header label for column B (row 1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "first_release_value"

# Apply the same style as A2 (date column, s=2) to the new date cells A3:A22
# before filling values, so row heights/format are consistent with the diff.
$ws.Range("A2").Copy()
$ws.Range("A3:A22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Column A: date serials (one row shifted down, new final row appended)
$ws.Range("A2").Value = 38717
$ws.Range("A3").Value = 39082
$ws.Range("A4").Value = 39447
$ws.Range("A5").Value = 39813
$ws.Range("A6").Value = 40178
$ws.Range("A7").Value = 40543
$ws.Range("A8").Value = 40908
$ws.Range("A9").Value = 41274
$ws.Range("A10").Value = 41639
$ws.Range("A11").Value = 42004
$ws.Range("A12").Value = 42369
$ws.Range("A13").Value = 42735
$ws.Range("A14").Value = 43100
$ws.Range("A15").Value = 43465
$ws.Range("A16").Value = 43830
$ws.Range("A17").Value = 44196
$ws.Range("A18").Value = 44561
$ws.Range("A19").Value = 44926
$ws.Range("A20").Value = 45291
$ws.Range("A21").Value = 45657
$ws.Range("A22").Value = 46022

# Column B: values shift down one row (old B2 becomes B3); B2 and the new
# final row (B22) are left empty.
$ws.Range("B2").ClearContents()
$ws.Range("B3").Value = 5.361718827437545
$ws.Range("B4").Value = 2.685523658736089
$ws.Range("B5").Value = 0.6019580713040096
$ws.Range("B6").Value = -2.848383350681438
$ws.Range("B7").Value = 0.5670099411379192
$ws.Range("B8").Value = 5.44843673085138
$ws.Range("B9").Value = -0.6152915357131694
$ws.Range("B10").Value = 0.4729742736614195
$ws.Range("B11").Value = 2.566421764830462
$ws.Range("B12").Value = 0.8407878010570302
$ws.Range("B13").Value = 2.194841098049016
$ws.Range("B14").Value = 2.643540836453884
$ws.Range("B15").Value = 2.959935600123309
$ws.Range("B16").Value = 3.452860220335019
$ws.Range("B17").Value = 1.607096457785584
$ws.Range("B18").Value = -3.136134057684858
$ws.Range("B19").Value = 1.613985729693268
$ws.Range("B20").Value = -1.992466799383086
$ws.Range("B21").Value = -3.350381746968589
$ws.Range("B22").ClearContents()
